$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits right
#    after the "IMAGEN:" run (it moves to the title paragraph below).
# ------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ------------------------------------------------------------------
# 2) Locate the paragraph that holds "SISTEMA DE NUMERACION CHINO"
#    (the recurso title, right after "* Titulo del recurso ..."),
#    being careful to use the FIRST occurrence in the document, not
#    the later, unrelated one further down.
# ------------------------------------------------------------------
$targetIndex = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx++
    if ($p.Range.Text -eq "SISTEMA DE NUMERACION CHINO`r") {
        $targetIndex = $idx
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the 'SISTEMA DE NUMERACION CHINO' paragraph"
}

$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range

# Split the paragraph in two: an empty paragraph first, then a new
# paragraph (inheriting the same paragraph/run formatting) that will
# hold the retyped title text.
$r.InsertParagraphBefore()

$p2 = $d.Paragraphs.Item($targetIndex + 1)
$full = $p2.Range
$target = $d.Range($full.Start, $full.End - 1)

# ------------------------------------------------------------------
# 3) Rewrite the paragraph's text as three runs - "S" / "istema de
#    numeracion" / " chino" - with the "_GoBack" bookmark sitting
#    between run 2 and run 3, matching the retyped text:
#    "Sistema de numeración chino". Inserting literal OOXML keeps
#    the runs distinct (rather than Word's usual same-format run
#    coalescing) and lets us place the bookmark exactly.
# ------------------------------------------------------------------
$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr>'
$oAcute = [char]0x00F3

$run1 = "<w:r>$rPr<w:t>S</w:t></w:r>"
$run2 = "<w:r>$rPr<w:t>istema de numeraci" + $oAcute + "n</w:t></w:r>"
$bookmark = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$run3 = "<w:r>$rPr<w:t xml:space=`"preserve`"> chino</w:t></w:r>"

$bodyXml = "<w:p>$run1$run2$bookmark$run3</w:p>"

$package = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  "<w:body>$bodyXml</w:body>" +
  '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($package)
